# Workbook / worksheet references
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1. Clear the leftover formatted-but-empty cells in C11:D14
#    (their content AND formatting are removed entirely)
# ---------------------------------------------------------------
$ws.Range("C11:D14").Clear()

# ---------------------------------------------------------------
# 2. Populate the new results table (rows 16-19) first so that the
#    shared-string table picks up the new text values in the same
#    order the original authoring session used.
# ---------------------------------------------------------------

# data_source for every new row
$ws.Range("B16").Value = "Mary Robert"

# model names
$ws.Range("A16").Value = "Linear Regression"
$ws.Range("A17").Value = "Lasso"
$ws.Range("A18").Value = "Ridge"
$ws.Range("A19").Value = "ElasticNet"

# remaining data_source cells (reuse the "Mary Robert" string)
$ws.Range("B17").Value = "Mary Robert"
$ws.Range("B18").Value = "Mary Robert"
$ws.Range("B19").Value = "Mary Robert"

# ---------------------------------------------------------------
# 3. Row 15 becomes a second header row for the new results table.
#    "test_MSE" then "train_MSE" are introduced in that order.
# ---------------------------------------------------------------
$ws.Range("F15").Value = "test_MSE"
$ws.Range("D15").Value = "train_MSE"

$ws.Range("A15").Value = "model"
$ws.Range("B15").Value = "data_source"
$ws.Range("C15").Value = "train_score"
$ws.Range("E15").Value = "test_score"

$ws.Range("C15").NumberFormat = "0.000"
$ws.Range("D15").NumberFormat = "0.000"
$ws.Range("E15").NumberFormat = "0.000"
$ws.Range("F15").NumberFormat = ".00"

# ---------------------------------------------------------------
# 4. Numeric values + number formats for the new data rows
# ---------------------------------------------------------------

# Row 16 - Linear Regression
$ws.Range("C16").Value = 0.02
$ws.Range("D16").Value = 0.98
$ws.Range("E16").Value = 0.01
$ws.Range("F16").Value = 1.38
$ws.Range("C16:F16").NumberFormat = ".00"

# Row 17 - Lasso
$ws.Range("C17").NumberFormat = ".00"
$ws.Range("E17").Value = 0.01
$ws.Range("F17").Value = 1.38
$ws.Range("E17:F17").NumberFormat = ".00"

# Row 18 - Ridge
$ws.Range("C18").NumberFormat = ".00"
$ws.Range("E18").Value = 0.01
$ws.Range("F18").Value = 1.38
$ws.Range("E18:F18").NumberFormat = ".00"

# Row 19 - ElasticNet
$ws.Range("C19").NumberFormat = ".00"
$ws.Range("E19").Value = 0.01
$ws.Range("F19").Value = 1.38
$ws.Range("E19:F19").NumberFormat = ".00"

# ---------------------------------------------------------------
# 5. Row 20 - trailing formatted blank cell
# ---------------------------------------------------------------
$ws.Range("E20").NumberFormat = ".00"

# ---------------------------------------------------------------
# 6. Re-fit column F now that it holds the new "test_MSE" header
#    and larger numeric values (matches Excel's recalculated
#    "best fit" width for the new content)
# ---------------------------------------------------------------
$ws.Columns.Item(6).ColumnWidth = 8.25

# ---------------------------------------------------------------
# 7. Leave the selection where the user's editing session ended
# ---------------------------------------------------------------
$ws.Range("C22").Select() | Out-Null
